# Populate the worksheet with location/city code data, matching the
# order strings are first introduced so the shared-string table lines
# up with the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows (2-4) first, column order A, B, C, D -----------------
$ws.Range("A2").Value = "NAEK176"
$ws.Range("B2").Value = "안중오거리"
$ws.Range("C2").Value = 31070
$ws.Range("D2").Value = "평택시"

$ws.Range("A3").Value = "NAEK177"
$ws.Range("B3").Value = "안중"
$ws.Range("C3").Value = 31070
$ws.Range("D3").Value = "평택시"

$ws.Range("A4").Value = "NAEK339"
$ws.Range("B4").Value = "아산온양"
$ws.Range("C4").Value = 34040
$ws.Range("D4").Value = "아산시"

# --- Header row (1) added afterwards, in B, C, A order ---------------
$ws.Range("B1").Value = "locationName"
$ws.Range("C1").Value = "cityCode"
$ws.Range("A1").Value = "locationCode"

# --- Styling: data rows get an explicit black font color -------------
$ws.Range("A2:D4").Font.Color = 0

# --- Selection ends on F8, matching the saved sheetView --------------
$null = $ws.Range("F8").Select()
